$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 2616.875
$ws.Range("J70").Value = 3142.5
$ws.Range("L70").Value = 9427.5
$ws.Range("N70").Value = -9967.5
$ws.Range("H73").Value = 2616.875
$ws.Range("J73").Value = 3142.5
$ws.Range("L73").Value = 9427.5
$ws.Range("N73").Value = -11299.5
$ws.Range("H137").Value = 40001804
$ws.Range("I137").Value = 62501068
$ws.Range("J137").Value = 3110.5557
$ws.Range("K137").Value = 187503204
$ws.Range("L137").Value = 9331.667099999999
$ws.Range("M137").Value = -187500654
$ws.Range("N137").Value = -14431.6671
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1033.9429
$ws.Range("I45").Value = 1008.1739
$ws.Range("J45").Value = 1083.3334
$ws.Range("K45").Value = 1008.1739
$ws.Range("L45").Value = 1083.3334
$ws.Range("M45").Value = -631.1739
$ws.Range("N45").Value = -1837.3334
$ws.Range("H61").Value = 2923.2068
$ws.Range("I61").Value = 1587
$ws.Range("J61").Value = 5109.727
$ws.Range("K61").Value = 1587
$ws.Range("L61").Value = 5109.727
$ws.Range("M61").Value = -1375
$ws.Range("N61").Value = -5533.727
$ws.Range("H74").Value = 8630.352999999999
$ws.Range("I74").Value = 1681.0667
$ws.Range("J74").Value = 60750
$ws.Range("K74").Value = 1681.0667
$ws.Range("L74").Value = 60750
$ws.Range("M74").Value = -807.0667000000001
$ws.Range("N74").Value = -62498
$ws.Range("H77").Value = 8630.352999999999
$ws.Range("I77").Value = 1681.0667
$ws.Range("J77").Value = 60750
$ws.Range("K77").Value = 8405.333500000001
$ws.Range("L77").Value = 303750
$ws.Range("M77").Value = -4037.333500000001
$ws.Range("N77").Value = -312486
$ws.Range("H132").Value = 3261.2222
$ws.Range("I132").Value = 2782.889
$ws.Range("J132").Value = 3739.5557
$ws.Range("K132").Value = 8348.667000000001
$ws.Range("L132").Value = 11218.6671
$ws.Range("M132").Value = -5818.667000000001
$ws.Range("N132").Value = -16278.6671
$ws.Range("H136").Value = 2923.2068
$ws.Range("I136").Value = 1587
$ws.Range("J136").Value = 5109.727
$ws.Range("K136").Value = 4761
$ws.Range("L136").Value = 15329.181
$ws.Range("M136").Value = -2211
$ws.Range("N136").Value = -20429.181
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3372.7273
$ws.Range("I134").Value = 2008.2273
$ws.Range("J134").Value = 6101.727
$ws.Range("K134").Value = 6024.6819
$ws.Range("L134").Value = 18305.181
$ws.Range("M134").Value = -3489.6819
$ws.Range("N134").Value = -23375.181
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2466.6667
$ws.Range("I31").Value = 2200
$ws.Range("K31").Value = 2200
$ws.Range("M31").Value = -1905
$ws.Range("H34").Value = 2466.6667
$ws.Range("I34").Value = 2200
$ws.Range("K34").Value = 2200
$ws.Range("M34").Value = -1998
$ws.Range("H58").Value = 2596.647
$ws.Range("I58").Value = 1655.3889
$ws.Range("J58").Value = 3655.5625
$ws.Range("K58").Value = 1655.3889
$ws.Range("L58").Value = 3655.5625
$ws.Range("M58").Value = -1452.3889
$ws.Range("N58").Value = -4061.5625
$ws.Range("H132").Value = 3664.8
$ws.Range("I132").Value = 2914
$ws.Range("K132").Value = 8742
$ws.Range("M132").Value = -6212
$ws.Range("H134").Value = 2170.853
$ws.Range("I134").Value = 1000
$ws.Range("J134").Value = 4317.4165
$ws.Range("K134").Value = 3000
$ws.Range("L134").Value = 12952.2495
$ws.Range("M134").Value = -465
$ws.Range("N134").Value = -18022.2495
$ws.Range("H136").Value = 2596.647
$ws.Range("I136").Value = 1655.3889
$ws.Range("J136").Value = 3655.5625
$ws.Range("K136").Value = 4966.1667
$ws.Range("L136").Value = 10966.6875
$ws.Range("M136").Value = -2416.1667
$ws.Range("N136").Value = -16066.6875
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 852.2
$ws.Range("I5").Value = 628.8461
$ws.Range("J5").Value = 1094.1666
$ws.Range("K5").Value = 1886.5383
$ws.Range("L5").Value = 3282.4998
$ws.Range("M5").Value = -1774.5383
$ws.Range("N5").Value = -3506.4998
$ws.Range("H80").Value = 1087.5
$ws.Range("I80").Value = 900
$ws.Range("J80").Value = 1150
$ws.Range("K80").Value = 2700
$ws.Range("L80").Value = 3450
$ws.Range("M80").Value = -1764
$ws.Range("N80").Value = -5322
$ws.Range("H83").Value = 1087.5
$ws.Range("I83").Value = 900
$ws.Range("J83").Value = 1150
$ws.Range("K83").Value = 8100
$ws.Range("L83").Value = 10350
$ws.Range("M83").Value = -3420
$ws.Range("N83").Value = -19710
$ws.Range("H114").Value = 1237.125
$ws.Range("I114").Value = 658.3
$ws.Range("J114").Value = 1500.2273
$ws.Range("K114").Value = 1974.9
$ws.Range("L114").Value = 4500.6819
$ws.Range("M114").Value = 1279.1
$ws.Range("N114").Value = -11008.6819
$ws.Range("J117").Value = 458.57144
$ws.Range("L117").Value = 1375.71432
$ws.Range("N117").Value = -8259.714319999999
$ws.Range("H135").Value = 852.2
$ws.Range("I135").Value = 628.8461
$ws.Range("J135").Value = 1094.1666
$ws.Range("K135").Value = 5659.6149
$ws.Range("L135").Value = 9847.499400000001
$ws.Range("M135").Value = -3124.6149
$ws.Range("N135").Value = -14917.4994
$ws.Range("H137").Value = 6315962
$ws.Range("I137").Value = 8335863.5
$ws.Range("K137").Value = 25007590.5
$ws.Range("M137").Value = -25002490.5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1387.2106
$ws.Range("I102").Value = 1136.6364
$ws.Range("K102").Value = 1136.6364
$ws.Range("M102").Value = 485.3635999999999
$ws.Range("H132").Value = 3084.2559
$ws.Range("I132").Value = 2861
$ws.Range("J132").Value = 3341
$ws.Range("K132").Value = 8583
$ws.Range("L132").Value = 10023
$ws.Range("M132").Value = -6053
$ws.Range("N132").Value = -15083
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1133.8334
$ws.Range("I46").Value = 1075.25
$ws.Range("K46").Value = 1075.25
$ws.Range("M46").Value = -887.25
$ws.Range("H132").Value = 3738.1
$ws.Range("I132").Value = 2897.32
$ws.Range("J132").Value = 5139.4
$ws.Range("K132").Value = 8691.960000000001
$ws.Range("L132").Value = 15418.2
$ws.Range("M132").Value = -6161.960000000001
$ws.Range("N132").Value = -20478.2
$ws.Range("H136").Value = 5096.607
$ws.Range("I136").Value = 2728.9375
$ws.Range("J136").Value = 8253.5
$ws.Range("K136").Value = 8186.8125
$ws.Range("L136").Value = 24760.5
$ws.Range("M136").Value = -5636.8125
$ws.Range("N136").Value = -29860.5
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H101").Value = 0
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()
$ws.Range("H132").Value = 20003648
$ws.Range("I132").Value = 31253154
$ws.Range("J132").Value = 4526.6665
$ws.Range("K132").Value = 93759462
$ws.Range("L132").Value = 13579.9995
$ws.Range("M132").Value = -93756932
$ws.Range("N132").Value = -18639.9995
$ws.Range("H136").Value = 11530804
$ws.Range("I136").Value = 27862474
$ws.Range("J136").Value = 2565.0588
$ws.Range("K136").Value = 83587422
$ws.Range("L136").Value = 7695.176399999999
$ws.Range("M136").Value = -83584872
$ws.Range("N136").Value = -12795.1764
$ws.Range("H138").Value = 50000
$ws.Range("J138").Value = 50000
$ws.Range("L138").Value = 50000
$ws.Range("N138").Value = -60280
